$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.448.01'
$ws.Range('E2').Value = '  +0.73%  '
$ws.Range('D3').Value = '1.618.23'
$ws.Range('E3').Value = '  +1.62%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.89'
$ws.Range('E5').Value = '  -0.10%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.498'
$ws.Range('E6').Value = '  -0.50%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('E9').Value = '  +0.25%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.22'
$ws.Range('E10').Value = '  +1.27%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0847'
$ws.Range('E11').Value = '  -0.56%  '
$ws.Range('D12').Value = '1.846.52'
$ws.Range('E12').Value = '  +1.63%  '
$ws.Range('D13').Value = '1.620.23'
$ws.Range('E13').Value = '  +1.56%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.02'
$ws.Range('E14').Value = '  +0.53%  '
$ws.Range('E15').Value = '  +0.24%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.84'
$ws.Range('E16').Value = '  +0.00%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '237.98'
$ws.Range('E17').Value = '  +10.35%  '
$ws.Range('D18').Value = '26.465.47'
$ws.Range('E18').Value = '  +0.80%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.77'
$ws.Range('E19').Value = '  +5.48%  '
$ws.Range('D20').Value = '0.0₃0725'
$ws.Range('E20').Value = '  +0.31%  '
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range('B23').Value = 'Avalanche'
$ws.Range('C23').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.10'
$ws.Range('E23').Value = '  +0.90%  '
$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.18'
$ws.Range('E24').Value = '  +4.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.07'
$ws.Range('E25').Value = '  +1.25%  '
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('E27').Value = '  +1.20%  '
$ws.Range('E28').Value = '  +0.31%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.55'
$ws.Range('E29').Value = '  +2.92%  '
$ws.Range('E30').Value = '  +0.62%  '
$ws.Range('E31').Value = '  +0.05%  '
$ws.Range('D32').Value = '1.529.53'
$ws.Range('E32').Value = '  +7.64%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.24'
$ws.Range('E33').Value = '  +1.41%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.97'
$ws.Range('E34').Value = '  +0.37%  '
$ws.Range('E35').Value = '  +6.71%  '
$ws.Range('E36').Value = '  -0.13%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.569'
$ws.Range('E37').Value = '  -0.27%  '
$ws.Range('E38').Value = '  +0.24%  '
$ws.Range('E39').Value = '  +0.65%  '
$ws.Range('E40').Value = '  +2.39%  '
$ws.Range('E42').Value = '  +1.85%  '
$ws.Range('D43').Value = '1.757.95'
$ws.Range('E43').Value = '  +1.67%  '
$ws.Range('E44').Value = '  +0.03%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.916'
$ws.Range('E45').Value = '  -1.59%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '61.57'
$ws.Range('E46').Value = '  +0.98%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '90.41'
$ws.Range('E47').Value = '  +4.22%  '
$ws.Range('E48').Value = '  +1.66%  '
$ws.Range('E49').Value = '  +0.23%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0962'
$ws.Range('E50').Value = '  +1.05%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.53'
$ws.Range('E51').Value = '  +1.44%  '
